$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (target stored width 42.140625 characters; the closest
# value this host's pixel-snapped ColumnWidth setter can reach is 42.1667)
$ws.Columns("A").ColumnWidth = 41.33

# Update the ID values in column A to the new 10-digit identifiers
$ws.Range("A2").Value = 1000000001
$ws.Range("A3").Value = 1000000002
$ws.Range("A4").Value = 1000000003
$ws.Range("A5").Value = 1000000004

# Move the selection to D2 (single cell) as in the saved view state
$ws.Range("D2").Select()
